$wb = $excel.ActiveWorkbook

# --- Update the "Status" text for the second data row (49c2e7af... file) ---
# This shared string is used by Overview!E3/F3, zh-cn!C3 and de-de!C3.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"

# --- zh-cn sheet: widen the "Error Detail" column (P) and record the handback error ---
# (ColumnWidth 39.17 round-trips to a stored sheet column width of exactly 40,
#  matching how this workbook already represents its other 40-wide columns.)
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsZhCn.Range("P3").Value = "Handback file name: dusyu1eb.ne3 is different with handoff file name: 49c2e7af-890e-4a4c-ad79-5631a57068cf.3f9ee1222676c4e760c567f148f6bb18ac882e9a.zh-cn."

# --- de-de sheet: widen the "Error Detail" column (P) and record the handback error ---
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Range("P3").Value = "Handback file name: dusyu1eb.ne3 is different with handoff file name: 49c2e7af-890e-4a4c-ad79-5631a57068cf.3f9ee1222676c4e760c567f148f6bb18ac882e9a.de-de."
